# Add cantrals by cantons
# Restructure the GE 2016 water cantons sheet:
#  - remove the old "units" row (row 2: Hiver/Hiver/Ete/Ete/Annee)
#  - rebuild the header row with the new column layout
#    (idx, idx2, Name, Date Start, Date End, (m3/s), (MW1), (MW2),
#     (GWh) Winter, (GWh) Summer, (GWh) Year)
#  - shift the data rows up accordingly (handled automatically by the
#    row delete) and update the selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old units row (previously row 2); this shifts the five
# data rows up by one and keeps all of their values/styles intact.
$ws.Rows(2).Delete()

# Rebuild row 1 with the new headers.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Give the F1:K1 headers their own style (Arial 9, general format),
# distinct from the plain A1:E1 headers.
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.ColorIndex = 1

# Match the saved selection/active cell.
$null = $ws.Range("A2:K2").Select()
